# Fruta / hortaliza, semanal
# Insert 3 new weekly price records (rows 261-263) into the Choclo sheet.
# This pushes the existing rows 261-312 down to 264-315 (dimension grows
# from A1:R312 to A1:R315) while the now-trailing rows 313-315 retain the
# content that used to sit at 310-312.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 261, shifting everything
# below (rows 261-312) down by three rows in one shot.
$ws.Range("A261:A263").EntireRow.Insert()

# --- New row 261 ---
$ws.Range("A261").Value = 9
$ws.Range("B261").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C261").Value = "Metropolitana"
$ws.Range("D261").Value = 44511
$ws.Range("E261").Value = 13
$ws.Range("F261").Value = 100112024
$ws.Range("G261").Value = "Choclo"
$ws.Range("H261").Value = "Choclero"
$ws.Range("I261").Value = "Primera"
$ws.Range("J261").Value = 52
$ws.Range("K261").Value = 28000
$ws.Range("L261").Value = 30000
$ws.Range("M261").Value = 29000
$ws.Range("N261").Value = "$/malla 50 unidades"
$ws.Range("O261").Value = "Región de Arica y Parinacota"
$ws.Range("P261").Value = 580
$ws.Range("Q261").Value = 50
$ws.Range("R261").Value = "Hortaliza"

# --- New row 262 ---
$ws.Range("A262").Value = 9
$ws.Range("B262").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C262").Value = "Metropolitana"
$ws.Range("D262").Value = 44511
$ws.Range("E262").Value = 13
$ws.Range("F262").Value = 100112024
$ws.Range("G262").Value = "Choclo"
$ws.Range("H262").Value = "Dulce o Americano"
$ws.Range("I262").Value = "Primera"
$ws.Range("J262").Value = 36
$ws.Range("K262").Value = 25000
$ws.Range("L262").Value = 26000
$ws.Range("M262").Value = 25500
$ws.Range("N262").Value = "$/caja 50 unidades"
$ws.Range("O262").Value = "Argentina"
$ws.Range("P262").Value = 510
$ws.Range("Q262").Value = 50
$ws.Range("R262").Value = "Hortaliza"

# --- New row 263 ---
$ws.Range("A263").Value = 9
$ws.Range("B263").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C263").Value = "Metropolitana"
$ws.Range("D263").Value = 44511
$ws.Range("E263").Value = 13
$ws.Range("F263").Value = 100112024
$ws.Range("G263").Value = "Choclo"
$ws.Range("H263").Value = "Dulce o Americano"
$ws.Range("I263").Value = "Primera"
$ws.Range("J263").Value = 79
$ws.Range("K263").Value = 25000
$ws.Range("L263").Value = 28000
$ws.Range("M263").Value = 26481
$ws.Range("N263").Value = "$/malla 70 unidades"
$ws.Range("O263").Value = "Región de Arica y Parinacota"
$ws.Range("P263").Value = 378
$ws.Range("Q263").Value = 70
$ws.Range("R263").Value = "Hortaliza"
